# Update column F (dSF) values on the active sheet to reflect the
# repulled data / recalculated means described in the commit message.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 7
    3  = -4
    4  = 3
    6  = -1
    9  = -3
    10 = 1
    11 = -7
    13 = 1
    14 = -4
    17 = -3
    18 = -1
    19 = 2
    20 = -4
    21 = 1
    22 = 4
    23 = -2
    24 = -3
    25 = -3
    27 = 3
    28 = 1
    29 = 2
    30 = -7
    32 = -5
    33 = -2
    34 = -3
    35 = 1
    37 = 2
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
